$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Step 1: Move the "Play Age of the Gods..." (bold) paragraph from the
# bottom of the document to right after the title (Heading1) paragraph,
# turning it into the new "Meta description" paragraph: the bold run's
# text becomes the label "Meta description" and the rest of the old
# subtitle sentence is appended right after it as a second, unbolded
# run - producing a "Meta description: <old subtitle>" paragraph.
# -----------------------------------------------------------------

# Find the bold "Play Age of the Gods Slots Free: Review & Guide"
# paragraph near the end of the document. (Range.Text includes the
# trailing paragraph-mark character, so strip that before comparing.
# Skip the first paragraph, which is the identically worded title.)
$oldLabel = "Play Age of the Gods Slots Free: Review & Guide"
$boldParaIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($i -gt 1 -and $t -eq $oldLabel) {
        $boldParaIndex = $i
    }
}
$boldPara = $d.Paragraphs($boldParaIndex)

# Cut the whole paragraph (including its end-of-paragraph mark) so the
# exact same run/formatting structure (bold run) is relocated rather
# than rebuilt from scratch.
$boldPara.Range.Cut()

# Paste it immediately after the title paragraph (paragraph 1), making
# it the new paragraph 2.
$titlePara = $d.Paragraphs(1)
$insertPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)
$insertPoint.Paste()

# Relabel the bold run's visible text ("Play Age of the Gods Slots
# Free: Review & Guide") to "Meta description", keeping it bold.
$newPara = $d.Paragraphs(2)
$labelStart = $newPara.Range.Start
$labelEnd = $labelStart + $oldLabel.Length
$labelRange = $d.Range($labelStart, $labelEnd)
$labelRange.Text = "Meta description"

# Build the continuation text (": Explore the immersive...") in a
# scratch spot borrowed from an existing plain (non-bold, non-italic)
# paragraph's END -- not its start, so as to not disturb any leading
# empty run the host paragraph may have -- then cut it from there and
# paste it right after the bold label. This guarantees the appended
# run picks up plain/default character formatting instead of
# inheriting Bold from "Meta description".
$hostPara = $d.Paragraphs(4)
$hostParaEnd = $hostPara.Range.End
$hostPoint = $d.Range($hostParaEnd - 1, $hostParaEnd - 1)
$appendText = ": Explore the immersive and mythical Age of Gods online slot, win four progressive jackpots, free spins, and expanding wilds. Play now on desktop and mobile."
$hostPoint.InsertAfter($appendText)
$hostRange = $d.Range($hostParaEnd - 1, $hostParaEnd - 1 + $appendText.Length)
$hostRange.Cut()

$newPara = $d.Paragraphs(2)
$newParaEnd = $newPara.Range.End
$pastePoint = $d.Range($newParaEnd - 1, $newParaEnd - 1)
$pastePoint.Paste()

# -----------------------------------------------------------------
# Step 2: Update the remaining subtitle paragraph (now italic, at the
# very end of the document) with the new image-prompt text.
# -----------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$subtitlePara = $d.Paragraphs($lastIndex)
$subStart = $subtitlePara.Range.Start
$subEnd = $subtitlePara.Range.End
$subTextRange = $d.Range($subStart, $subEnd - 1)
$subTextRange.Text = 'Create a feature image for the game "Age of the Gods" that highlights the Greek mythology theme and features a happy Maya warrior with glasses. The image should be in cartoon style and should include Mount Olympus and the main deities such as Athena, Zeus, Poseidon, Hades, and Hercules. The Maya warrior can be positioned in the center of the image, holding a slot machine lever or spinning a wheel with a confident and excited expression. The overall tone should be vibrant and colorful, highlighting the adventurous and rewarding nature of this popular online slot game.'
